$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2301.2
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 2753
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 2753
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -4001
$ws.Range("H64").Value = 2720.923
$ws.Range("J64").Value = 2846.25
$ws.Range("L64").Value = 2846.25
$ws.Range("N64").Value = -3342.25
$ws.Range("H65").Value = 2301.2
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 2753
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 13765
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -20005
$ws.Range("H67").Value = 2720.923
$ws.Range("J67").Value = 2846.25
$ws.Range("L67").Value = 2846.25
$ws.Range("N67").Value = -4562.25
$ws.Range("H76").Value = 93490.73
$ws.Range("I76").Value = 93490.73
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 93490.73
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -93175.73
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 93490.73
$ws.Range("I79").Value = 93490.73
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 93490.73
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -92398.73
$ws.Range("N79").ClearContents()
$ws.Range("H137").Value = 1376.6072
$ws.Range("I137").Value = 987.9048
$ws.Range("J137").Value = 2542.7144
$ws.Range("K137").Value = 2963.7144
$ws.Range("L137").Value = 7628.1432
$ws.Range("M137").Value = -413.7143999999998
$ws.Range("N137").Value = -12728.1432
$ws.Range("H141").Value = 1928.1
$ws.Range("I141").Value = 1392.4722
$ws.Range("J141").Value = 6748.75
$ws.Range("K141").Value = 4177.4166
$ws.Range("L141").Value = 20246.25
$ws.Range("M141").Value = 1002.5834
$ws.Range("N141").Value = -30606.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17054.361
$ws.Range("I32").Value = 16658.438
$ws.Range("J32").Value = 22333.334
$ws.Range("K32").Value = 16658.438
$ws.Range("L32").Value = 22333.334
$ws.Range("M32").Value = -16371.438
$ws.Range("N32").Value = -22907.334

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 990.2
$ws.Range("I16").Value = 960
$ws.Range("K16").Value = 960
$ws.Range("M16").Value = -673
$ws.Range("H31").Value = 3033538.8
$ws.Range("I31").Value = 2481.2104
$ws.Range("J31").Value = 9808844
$ws.Range("K31").Value = 2481.2104
$ws.Range("L31").Value = 9808844
$ws.Range("M31").Value = -2186.2104
$ws.Range("N31").Value = -9809434
$ws.Range("H34").Value = 3033538.8
$ws.Range("I34").Value = 2481.2104
$ws.Range("J34").Value = 9808844
$ws.Range("K34").Value = 2481.2104
$ws.Range("L34").Value = 9808844
$ws.Range("M34").Value = -2279.2104
$ws.Range("N34").Value = -9809248
$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50588
$ws.Range("H58").Value = 1081.1364
$ws.Range("I58").Value = 1099.5238
$ws.Range("J58").Value = 695
$ws.Range("K58").Value = 1099.5238
$ws.Range("L58").Value = 695
$ws.Range("M58").Value = -896.5237999999999
$ws.Range("N58").Value = -1101
$ws.Range("H64").Value = 42499.75
$ws.Range("I64").Value = 40000
$ws.Range("J64").Value = 44999.5
$ws.Range("K64").Value = 40000
$ws.Range("L64").Value = 44999.5
$ws.Range("M64").Value = -39752
$ws.Range("N64").Value = -45495.5
$ws.Range("H67").Value = 42499.75
$ws.Range("I67").Value = 40000
$ws.Range("J67").Value = 44999.5
$ws.Range("K67").Value = 40000
$ws.Range("L67").Value = 44999.5
$ws.Range("M67").Value = -39142
$ws.Range("N67").Value = -46715.5
$ws.Range("H99").Value = 2485.2
$ws.Range("I99").Value = 2142.2222
$ws.Range("J99").Value = 2999.6667
$ws.Range("K99").Value = 2142.2222
$ws.Range("L99").Value = 2999.6667
$ws.Range("M99").Value = -644.2222000000002
$ws.Range("N99").Value = -5995.6667
$ws.Range("H113").Value = 990.2
$ws.Range("I113").Value = 960
$ws.Range("K113").Value = 960
$ws.Range("M113").Value = 1210
$ws.Range("H126").Value = 2485.2
$ws.Range("I126").Value = 2142.2222
$ws.Range("J126").Value = 2999.6667
$ws.Range("K126").Value = 6426.6666
$ws.Range("L126").Value = 8999.000100000001
$ws.Range("M126").Value = -3956.6666
$ws.Range("N126").Value = -13939.0001
$ws.Range("H136").Value = 1081.1364
$ws.Range("I136").Value = 1099.5238
$ws.Range("J136").Value = 695
$ws.Range("K136").Value = 3298.5714
$ws.Range("L136").Value = 2085
$ws.Range("M136").Value = -748.5713999999998
$ws.Range("N136").Value = -7185

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 108.4
$ws.Range("I38").Value = 10.2
$ws.Range("J38").Value = 157.5
$ws.Range("K38").Value = 30.6
$ws.Range("L38").Value = 472.5
$ws.Range("M38").Value = 316.4
$ws.Range("N38").Value = -1166.5
$ws.Range("H62").Value = 3249.6667
$ws.Range("I62").Value = 2499.4285
$ws.Range("J62").Value = 4300
$ws.Range("K62").Value = 7498.2855
$ws.Range("L62").Value = 12900
$ws.Range("M62").Value = -6812.2855
$ws.Range("N62").Value = -14272
$ws.Range("H65").Value = 3249.6667
$ws.Range("I65").Value = 2499.4285
$ws.Range("J65").Value = 4300
$ws.Range("K65").Value = 22494.8565
$ws.Range("L65").Value = 38700
$ws.Range("M65").Value = -19062.8565
$ws.Range("N65").Value = -45564

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 62504390
$ws.Range("I70").Value = 93754280
$ws.Range("J70").Value = 4613.5
$ws.Range("K70").Value = 93754280
$ws.Range("L70").Value = 4613.5
$ws.Range("M70").Value = -93754010
$ws.Range("N70").Value = -5153.5
$ws.Range("H73").Value = 62504390
$ws.Range("I73").Value = 93754280
$ws.Range("J73").Value = 4613.5
$ws.Range("K73").Value = 93754280
$ws.Range("L73").Value = 4613.5
$ws.Range("M73").Value = -93753344
$ws.Range("N73").Value = -6485.5
$ws.Range("H80").Value = 5701.5
$ws.Range("I80").Value = 4659.2856
$ws.Range("J80").Value = 8133.3335
$ws.Range("K80").Value = 4659.2856
$ws.Range("L80").Value = 8133.3335
$ws.Range("M80").Value = -3661.2856
$ws.Range("N80").Value = -10129.3335
$ws.Range("H83").Value = 5701.5
$ws.Range("I83").Value = 4659.2856
$ws.Range("J83").Value = 8133.3335
$ws.Range("K83").Value = 23296.428
$ws.Range("L83").Value = 40666.6675
$ws.Range("M83").Value = -18304.428
$ws.Range("N83").Value = -50650.6675
$ws.Range("H122").Value = 27779800
$ws.Range("I122").Value = 38463416
$ws.Range("J122").Value = 2397.2
$ws.Range("K122").Value = 115390248
$ws.Range("L122").Value = 7191.599999999999
$ws.Range("M122").Value = -115387798
$ws.Range("N122").Value = -12091.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 455.3
$ws.Range("I113").Value = 526
$ws.Range("J113").Value = 384.6
$ws.Range("K113").Value = 1578
$ws.Range("L113").Value = 1153.8
$ws.Range("M113").Value = 592
$ws.Range("N113").Value = -5493.8
$ws.Range("H132").Value = 1619.102
$ws.Range("I132").Value = 1523.6945
$ws.Range("J132").Value = 1883.3077
$ws.Range("K132").Value = 4571.083500000001
$ws.Range("L132").Value = 5649.9231
$ws.Range("M132").Value = -2041.083500000001
$ws.Range("N132").Value = -10709.9231
$ws.Range("H136").Value = 941.6271400000001
$ws.Range("I136").Value = 904.4035
$ws.Range("J136").Value = 2002.5
$ws.Range("K136").Value = 2713.2105
$ws.Range("L136").Value = 6007.5
$ws.Range("M136").Value = -163.2105000000001
$ws.Range("N136").Value = -11107.5
